## Edit: move the "product structure / repository" lookup tables from
## Sheet1 (rows 31-60) into a brand-new Sheet2, update Sheet1's used range
## / view state, add the "Clay Tandoor" category to row 28 (column L),
## and add the "Clay Tandoor" entry as a new shared string.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Insert the new worksheet right after Sheet1 --------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# --- 2. Populate Sheet2 with the repository/lookup tables ---------------
# Block: Table1 (SizeID -> Size) / Table2 (ProductId -> SizeID) / Table1
# (SizeID, ShapeID -> WeightId)
$ws2.Cells.Item(3, 2).Value  = "Table1"
$ws2.Cells.Item(3, 5).Value  = "Table2"
$ws2.Cells.Item(3, 8).Value  = "Table1"

$ws2.Cells.Item(4, 2).Value  = "SizeID"
$ws2.Cells.Item(4, 3).Value  = "Size"
$ws2.Cells.Item(4, 5).Value  = "ProductId"
$ws2.Cells.Item(4, 6).Value  = "SizeID"
$ws2.Cells.Item(4, 8).Value  = "SizeID"
$ws2.Cells.Item(4, 9).Value  = "ShapeID"
$ws2.Cells.Item(4, 10).Value = "WeightId"

$ws2.Cells.Item(5, 2).Value  = 2
$ws2.Cells.Item(5, 3).Value  = 34
$ws2.Cells.Item(5, 5).Value  = 1
$ws2.Cells.Item(5, 6).Value  = 2
$ws2.Cells.Item(5, 8).Value  = 2
$ws2.Cells.Item(5, 9).Value  = 3
$ws2.Cells.Item(5, 10).Value = 4

$ws2.Cells.Item(6, 2).Value  = 3
$ws2.Cells.Item(6, 3).Value  = 32
$ws2.Cells.Item(6, 5).Value  = 1
$ws2.Cells.Item(6, 6).Value  = 3
$ws2.Cells.Item(6, 8).Value  = 3
$ws2.Cells.Item(6, 9).Value  = 3
$ws2.Cells.Item(6, 10).Value = 5

$ws2.Cells.Item(7, 2).Value  = 4
$ws2.Cells.Item(7, 3).Value  = 30
$ws2.Cells.Item(7, 5).Value  = 1
$ws2.Cells.Item(7, 6).Value  = 4
$ws2.Cells.Item(7, 8).Value  = 4
$ws2.Cells.Item(7, 9).Value  = 3
$ws2.Cells.Item(7, 10).Value = 6

$ws2.Cells.Item(8, 2).Value  = 5
$ws2.Cells.Item(8, 3).Value  = 18
$ws2.Cells.Item(8, 8).Value  = 5
$ws2.Cells.Item(8, 9).Value  = 4
$ws2.Cells.Item(8, 10).Value = 4

$ws2.Cells.Item(9, 2).Value  = 6
$ws2.Cells.Item(9, 3).Value  = 12
$ws2.Cells.Item(9, 8).Value  = 6
$ws2.Cells.Item(9, 9).Value  = 4
$ws2.Cells.Item(9, 10).Value = 5

# Table3: ShapeId, Shape, Weight
$ws2.Cells.Item(12, 2).Value = "Table3"
$ws2.Cells.Item(13, 2).Value = "ShapeId"
$ws2.Cells.Item(13, 3).Value = "Shape"
$ws2.Cells.Item(13, 4).Value = "Weight"
$ws2.Cells.Item(14, 2).Value = 3
$ws2.Cells.Item(14, 3).Value = "Round"
$ws2.Cells.Item(14, 4).Value = 50
$ws2.Cells.Item(15, 2).Value = 4
$ws2.Cells.Item(15, 3).Value = "Square"
$ws2.Cells.Item(15, 4).Value = 100
$ws2.Cells.Item(16, 2).Value = 2
$ws2.Cells.Item(16, 3).Value = "Baral"
$ws2.Cells.Item(16, 4).Value = 150
$ws2.Cells.Item(17, 2).Value = 2
$ws2.Cells.Item(17, 3).Value = "Square"

# Table4: WeightID, ShapeID, Weight
$ws2.Cells.Item(18, 2).Value = "Table4"
$ws2.Cells.Item(19, 2).Value = "WeightID"
$ws2.Cells.Item(19, 3).Value = "ShapeID"
$ws2.Cells.Item(19, 4).Value = "Weight"
$ws2.Cells.Item(20, 2).Value = 4
$ws2.Cells.Item(20, 3).Value = 3
$ws2.Cells.Item(20, 4).Value = 50
$ws2.Cells.Item(21, 4).Value = 100
$ws2.Cells.Item(22, 4).Value = 150

# Table5: ProductId
$ws2.Cells.Item(24, 2).Value = "Table5"
$ws2.Cells.Item(25, 2).Value = "ProductId"
$ws2.Cells.Item(26, 2).Value = 1

# Table6: ProductID, SizeID, ShapeID, WeightId
$ws2.Cells.Item(28, 2).Value = "Table6"
$ws2.Cells.Item(29, 2).Value = "ProductID"
$ws2.Cells.Item(29, 3).Value = "SizeID"
$ws2.Cells.Item(29, 4).Value = "ShapeID"
$ws2.Cells.Item(29, 5).Value = "WeightId"
$ws2.Cells.Item(30, 2).Value = 1
$ws2.Cells.Item(30, 3).Value = 1
$ws2.Cells.Item(30, 4).Value = 2
$ws2.Cells.Item(30, 5).Value = 4
$ws2.Cells.Item(31, 2).Value = 1
$ws2.Cells.Item(31, 3).Value = 2
$ws2.Cells.Item(31, 4).Value = 2
$ws2.Cells.Item(31, 5).Value = 5
$ws2.Cells.Item(32, 2).Value = 1
$ws2.Cells.Item(32, 3).Value = 3
$ws2.Cells.Item(32, 4).Value = 2
$ws2.Cells.Item(32, 5).Value = 6

# --- 3. Sheet2's view: select J4:J9, activate cell J4 --------------------
$ws2.Range("J4:J9").Select()
$ws2.Application.ActiveWindow.RangeSelection.Item(1).Activate()

# --- 4. Add the new "Clay Tandoor" category to Sheet1's product table ---
$ws1.Cells.Item(28, 12).Value = "Clay Tandoor"

# --- 5. Remove the now-relocated lookup tables from Sheet1 ---------------
$ws1.Range("A30:L60").EntireRow.Delete()

# --- 6. Restore Sheet1's view state (scrolled to top, selection at A3) --
$ws1.Application.ActiveWindow.ScrollRow = 2
$ws1.Range("A3").Select()

# --- 7. Make Sheet2 the active/visible tab -------------------------------
$ws2.Activate()
